$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2664.04
$ws.Range("I40").Value = 1740.8182
$ws.Range("J40").Value = 3389.4285
$ws.Range("K40").Value = 1740.8182
$ws.Range("L40").Value = 3389.4285
$ws.Range("M40").Value = -1565.8182
$ws.Range("N40").Value = -3739.4285
$ws.Range("H51").Value = 11112580
$ws.Range("I51").Value = 1400
$ws.Range("J51").Value = 12347155
$ws.Range("K51").Value = 1400
$ws.Range("L51").Value = 12347155
$ws.Range("M51").Value = -916
$ws.Range("N51").Value = -12348123
$ws.Range("H64").Value = 3164.6072
$ws.Range("I64").Value = 3112.5
$ws.Range("J64").Value = 3185.45
$ws.Range("K64").Value = 3112.5
$ws.Range("L64").Value = 3185.45
$ws.Range("M64").Value = -2864.5
$ws.Range("N64").Value = -3681.45
$ws.Range("H67").Value = 3164.6072
$ws.Range("I67").Value = 3112.5
$ws.Range("J67").Value = 3185.45
$ws.Range("K67").Value = 3112.5
$ws.Range("L67").Value = 3185.45
$ws.Range("M67").Value = -2254.5
$ws.Range("N67").Value = -4901.45
$ws.Range("H86").Value = 11914206
$ws.Range("I86").Value = 16678588
$ws.Range("J86").Value = 3253
$ws.Range("K86").Value = 16678588
$ws.Range("L86").Value = 3253
$ws.Range("M86").Value = -16677465
$ws.Range("N86").Value = -5499
$ws.Range("H89").Value = 11914206
$ws.Range("I89").Value = 16678588
$ws.Range("J89").Value = 3253
$ws.Range("K89").Value = 83392940
$ws.Range("L89").Value = 16265
$ws.Range("M89").Value = -83387324
$ws.Range("N89").Value = -27497
$ws.Range("H94").Value = 4076.25
$ws.Range("I94").Value = 4076.25
$ws.Range("K94").Value = 4076.25
$ws.Range("M94").Value = -3625.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10160.297
$ws.Range("I32").Value = 3085.9565
$ws.Range("J32").Value = 21782.428
$ws.Range("K32").Value = 3085.9565
$ws.Range("L32").Value = 21782.428
$ws.Range("M32").Value = -2798.9565
$ws.Range("N32").Value = -22356.428
$ws.Range("H76").Value = 21600
$ws.Range("I76").Value = 15000
$ws.Range("J76").Value = 26000
$ws.Range("K76").Value = 15000
$ws.Range("L76").Value = 26000
$ws.Range("M76").Value = -14662
$ws.Range("N76").Value = -26676
$ws.Range("H79").Value = 21600
$ws.Range("I79").Value = 15000
$ws.Range("J79").Value = 26000
$ws.Range("K79").Value = 15000
$ws.Range("L79").Value = 26000
$ws.Range("M79").Value = -13830
$ws.Range("N79").Value = -28340
$ws.Range("H92").Value = 20480
$ws.Range("J92").Value = 20480
$ws.Range("L92").Value = 20480
$ws.Range("N92").Value = -25472

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 10379.333
$ws.Range("I82").Value = 9073.429
$ws.Range("J82").Value = 14950
$ws.Range("K82").Value = 9073.429
$ws.Range("L82").Value = 14950
$ws.Range("M82").Value = -8690.429
$ws.Range("N82").Value = -15716
$ws.Range("H85").Value = 10379.333
$ws.Range("I85").Value = 9073.429
$ws.Range("J85").Value = 14950
$ws.Range("K85").Value = 9073.429
$ws.Range("L85").Value = 14950
$ws.Range("M85").Value = -7747.429
$ws.Range("N85").Value = -17602
$ws.Range("H86").Value = 5407022.5
$ws.Range("I86").Value = 6898143
$ws.Range("J86").Value = 1709.625
$ws.Range("K86").Value = 6898143
$ws.Range("L86").Value = 1709.625
$ws.Range("M86").Value = -6897020
$ws.Range("N86").Value = -3955.625
$ws.Range("H89").Value = 5407022.5
$ws.Range("I89").Value = 6898143
$ws.Range("J89").Value = 1709.625
$ws.Range("K89").Value = 34490715
$ws.Range("L89").Value = 8548.125
$ws.Range("M89").Value = -34485099
$ws.Range("N89").Value = -19780.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 48880
$ws.Range("J88").Value = 48880
$ws.Range("L88").Value = 48880
$ws.Range("N88").Value = -49692
$ws.Range("H91").Value = 48880
$ws.Range("J91").Value = 48880
$ws.Range("L91").Value = 48880
$ws.Range("N91").Value = -51688

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 672.43634
$ws.Range("J113").Value = 891.6799999999999
$ws.Range("L113").Value = 2675.04
$ws.Range("N113").Value = -7015.04

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1343.125
$ws.Range("I22").Value = 1716.6666
$ws.Range("J22").Value = 1119
$ws.Range("K22").Value = 1716.6666
$ws.Range("L22").Value = 1119
$ws.Range("M22").Value = -1421.6666
$ws.Range("N22").Value = -1709
$ws.Range("H27").Value = 1343.125
$ws.Range("I27").Value = 1716.6666
$ws.Range("J27").Value = 1119
$ws.Range("K27").Value = 1716.6666
$ws.Range("L27").Value = 1119
$ws.Range("M27").Value = -1609.6666
$ws.Range("N27").Value = -1333
$ws.Range("H46").Value = 1541.1818
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 1615.3
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 1615.3
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -1991.3
$ws.Range("H55").Value = 173.9
$ws.Range("I55").Value = 149.27272
$ws.Range("J55").Value = 204
$ws.Range("K55").Value = 149.27272
$ws.Range("L55").Value = 204
$ws.Range("M55").Value = 23.72728000000001
$ws.Range("N55").Value = -550
$ws.Range("H64").Value = 7037.5
$ws.Range("J64").Value = 7037.5
$ws.Range("L64").Value = 7037.5
$ws.Range("N64").Value = -7487.5
$ws.Range("H67").Value = 7037.5
$ws.Range("J67").Value = 7037.5
$ws.Range("L67").Value = 7037.5
$ws.Range("N67").Value = -8597.5
$ws.Range("H94").Value = 16875
$ws.Range("J94").Value = 16875
$ws.Range("L94").Value = 16875
$ws.Range("N94").Value = -18227

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 39362.5
$ws.Range("J82").Value = 39362.5
$ws.Range("L82").Value = 39362.5
$ws.Range("N82").Value = -40128.5
$ws.Range("H85").Value = 39362.5
$ws.Range("J85").Value = 39362.5
$ws.Range("L85").Value = 39362.5
$ws.Range("N85").Value = -42014.5
$ws.Range("H100").Value = 6992.6665
$ws.Range("I100").Value = 8974
$ws.Range("J100").Value = 1841.2
$ws.Range("K100").Value = 17948
$ws.Range("L100").Value = 3682.4
$ws.Range("M100").Value = -17407
$ws.Range("N100").Value = -4764.4
